# Updated symbol list on Wed Feb  1 22:44:51 UTC 2023 with GitHub Actions
#
# Refreshes the crypto ticker table (Coin / Link / Price / Volume(1h)) on
# Sheet1 with the latest pull from coinranking.com. A handful of rows also
# moved position in the source ranking, so B/C (Coin, Link) get rewritten
# alongside D/E (Price, Volume) for those rows.
#
# Price/Volume are stored as plain text in this sheet (not real numbers),
# so cells that look numeric ("316.29", "1.83%") need an explicit Text
# number format before the write - otherwise Excel helpfully "fixes" them
# into actual numbers/percentages, which is not what this sheet wants.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Col = 4; Val = '316.29'; Text = $true },
    @{ Row = 2; Col = 5; Val = '1.83%'; Text = $true },
    @{ Row = 3; Col = 4; Val = '37.82'; Text = $true },
    @{ Row = 3; Col = 5; Val = '1.42%'; Text = $true },
    @{ Row = 4; Col = 4; Val = '5.186'; Text = $true },
    @{ Row = 4; Col = 5; Val = '1.24%'; Text = $true },
    @{ Row = 5; Col = 4; Val = '0.07985'; Text = $true },
    @{ Row = 5; Col = 5; Val = '1.99%'; Text = $true },
    @{ Row = 6; Col = 2; Val = 'KuCoinToken'; Text = $false },
    @{ Row = 6; Col = 3; Val = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'; Text = $false },
    @{ Row = 6; Col = 4; Val = '8.543'; Text = $true },
    @{ Row = 6; Col = 5; Val = '3.49%'; Text = $true },
    @{ Row = 7; Col = 2; Val = 'FTXToken'; Text = $false },
    @{ Row = 7; Col = 3; Val = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'; Text = $false },
    @{ Row = 7; Col = 4; Val = '1.933'; Text = $true },
    @{ Row = 7; Col = 5; Val = '1.40%'; Text = $true },
    @{ Row = 8; Col = 2; Val = 'BTSEToken'; Text = $false },
    @{ Row = 8; Col = 3; Val = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'; Text = $false },
    @{ Row = 8; Col = 4; Val = '2.963'; Text = $true },
    @{ Row = 8; Col = 5; Val = '0.36%'; Text = $true },
    @{ Row = 9; Col = 2; Val = 'MXToken'; Text = $false },
    @{ Row = 9; Col = 3; Val = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; Text = $false },
    @{ Row = 9; Col = 4; Val = '0.9436'; Text = $true },
    @{ Row = 9; Col = 5; Val = '2.74%'; Text = $true },
    @{ Row = 10; Col = 2; Val = 'LiechtensteinCryptoassetsExchange'; Text = $false },
    @{ Row = 10; Col = 3; Val = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'; Text = $false },
    @{ Row = 10; Col = 4; Val = '0.1307'; Text = $true },
    @{ Row = 10; Col = 5; Val = '9.14%'; Text = $true },
    @{ Row = 11; Col = 2; Val = 'WazirX'; Text = $false },
    @{ Row = 11; Col = 3; Val = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'; Text = $false },
    @{ Row = 11; Col = 4; Val = '0.1937'; Text = $true },
    @{ Row = 11; Col = 5; Val = '1.19%'; Text = $true },
    @{ Row = 12; Col = 2; Val = 'MandalaExchangeToken'; Text = $false },
    @{ Row = 12; Col = 3; Val = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; Text = $false },
    @{ Row = 12; Col = 4; Val = '0.09083'; Text = $true },
    @{ Row = 12; Col = 5; Val = '1.01%'; Text = $true },
    @{ Row = 13; Col = 2; Val = 'BitrueCoin'; Text = $false },
    @{ Row = 13; Col = 3; Val = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; Text = $false },
    @{ Row = 13; Col = 4; Val = '0.03407'; Text = $true },
    @{ Row = 13; Col = 5; Val = '1.68%'; Text = $true },
    @{ Row = 14; Col = 2; Val = 'BitMartToken'; Text = $false },
    @{ Row = 14; Col = 3; Val = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; Text = $false },
    @{ Row = 14; Col = 4; Val = '0.09534'; Text = $true },
    @{ Row = 14; Col = 5; Val = '-0.63%'; Text = $true },
    @{ Row = 15; Col = 2; Val = 'BitForexToken'; Text = $false },
    @{ Row = 15; Col = 3; Val = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; Text = $false },
    @{ Row = 15; Col = 4; Val = '0.001390'; Text = $true },
    @{ Row = 15; Col = 5; Val = '0.93%'; Text = $true },
    @{ Row = 16; Col = 2; Val = 'TigerCash'; Text = $false },
    @{ Row = 16; Col = 3; Val = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'; Text = $false },
    @{ Row = 16; Col = 4; Val = '0.005879'; Text = $true },
    @{ Row = 16; Col = 5; Val = '2.68%'; Text = $true },
    @{ Row = 17; Col = 2; Val = 'LEO'; Text = $false },
    @{ Row = 17; Col = 3; Val = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; Text = $false },
    @{ Row = 17; Col = 4; Val = '3.434'; Text = $true },
    @{ Row = 17; Col = 5; Val = '-2.94%'; Text = $true },
    @{ Row = 18; Col = 2; Val = 'GateToken'; Text = $false },
    @{ Row = 18; Col = 3; Val = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'; Text = $false },
    @{ Row = 18; Col = 4; Val = '4.483'; Text = $true },
    @{ Row = 18; Col = 5; Val = '1.21%'; Text = $true },
    @{ Row = 19; Col = 5; Val = '2.16%'; Text = $true },
    @{ Row = 20; Col = 4; Val = '6.623'; Text = $true },
    @{ Row = 20; Col = 5; Val = '26.37%'; Text = $true },
    @{ Row = 21; Col = 4; Val = '0.1304'; Text = $true },
    @{ Row = 21; Col = 5; Val = '1.59%'; Text = $true },
    @{ Row = 22; Col = 5; Val = '-6.65%'; Text = $true },
    @{ Row = 23; Col = 4; Val = '0.04371'; Text = $true },
    @{ Row = 23; Col = 5; Val = '0.35%'; Text = $true },
    @{ Row = 24; Col = 4; Val = '0.001228'; Text = $true },
    @{ Row = 24; Col = 5; Val = '-1.79%'; Text = $true },
    @{ Row = 25; Col = 5; Val = '-8.54%'; Text = $true },
    @{ Row = 26; Col = 4; Val = '0.0001328'; Text = $true },
    @{ Row = 26; Col = 5; Val = '-2.53%'; Text = $true },
    @{ Row = 27; Col = 4; Val = '0.0003983'; Text = $true },
    @{ Row = 27; Col = 5; Val = '-0.32%'; Text = $true },
    @{ Row = 39; Col = 4; Val = '0.02405'; Text = $true },
    @{ Row = 39; Col = 5; Val = '6.40%'; Text = $true },
    @{ Row = 40; Col = 4; Val = '0.05153'; Text = $true },
    @{ Row = 40; Col = 5; Val = '2.13%'; Text = $true },
    @{ Row = 41; Col = 4; Val = '0.007615'; Text = $true },
    @{ Row = 41; Col = 5; Val = '1.91%'; Text = $true },
    @{ Row = 42; Col = 4; Val = '0.1400'; Text = $true },
    @{ Row = 42; Col = 5; Val = '3.67%'; Text = $true },
    @{ Row = 43; Col = 4; Val = '0.008571'; Text = $true },
    @{ Row = 43; Col = 5; Val = '-5.32%'; Text = $true },
    @{ Row = 44; Col = 4; Val = '0.002106'; Text = $true },
    @{ Row = 44; Col = 5; Val = '7.83%'; Text = $true },
    @{ Row = 45; Col = 4; Val = '0.008740'; Text = $true },
    @{ Row = 45; Col = 5; Val = '-6.10%'; Text = $true },
    @{ Row = 46; Col = 4; Val = '0.00006486'; Text = $true },
    @{ Row = 46; Col = 5; Val = '-1.21%'; Text = $true },
    @{ Row = 47; Col = 4; Val = '0.00000000749'; Text = $true },
    @{ Row = 47; Col = 5; Val = '-0.32%'; Text = $true },
    @{ Row = 48; Col = 4; Val = '0.002862'; Text = $true },
    @{ Row = 48; Col = 5; Val = '-15.14%'; Text = $true },
    @{ Row = 49; Col = 4; Val = '0.001687'; Text = $true },
    @{ Row = 49; Col = 5; Val = '68.40%'; Text = $true },
    @{ Row = 50; Col = 4; Val = '0.00002096'; Text = $true },
    @{ Row = 50; Col = 5; Val = '-0.32%'; Text = $true },
    @{ Row = 51; Col = 4; Val = '0.0001996'; Text = $true },
    @{ Row = 51; Col = 5; Val = '-0.32%'; Text = $true }
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u.Row, $u.Col)
    if ($u.Text) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Val
}

Write-Output "Applied $($updates.Count) cell updates"
